$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SwateTemplateMetadata")
$ws.Name = "isa_template"
$ws.Range("E13").Clear()
$ws.Range("D14").Clear()
